$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: locate the row number of a country by its exact name in column A
function Get-RowByCountry($ws, $name) {
    $searchRange = $ws.Range("A4:A205")
    $found = $searchRange.Find($name, $null, $null, 1, $null, $null, $null, $null, $null)
    return $found.Row
}

# New case numbers (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the countries whose figures changed
$updates = @{
    "Estados Unidos" = @{ B = 101159; C = 15724; D = 2465;  E = 97135; F = 2463; G = 264; H = 1559 }
    "China"           = @{ B = 81340;  C = 0;     D = 74588; E = 3460;  F = 1034; G = 0;   H = 3292 }
    "España"          = @{ B = 65719;  C = 7933;  D = 9357;  E = 51224; F = 4165; G = 773; H = 5138 }
    "Canada"          = @{ B = 4633;   C = 590;   D = 258;   E = 4322;  F = 120;  G = 14;  H = 53 }
    "Uganda"          = @{ B = 23;     C = 9;     D = 0;     E = 23;    F = 0;    G = 0;   H = 0 }
    "Bermudas"        = @{ B = 17;     C = 2;     D = 2;     E = 15;    F = 0;    G = 0;   H = 0 }
}

foreach ($country in $updates.Keys) {
    $row = Get-RowByCountry $ws $country
    $vals = $updates[$country]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# Update the "last updated" timestamp shown in the title row
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 23:28"

# Re-sort the country table (rows 4-205) descending by "Casos totales" (column B),
# which is how the sheet is normally kept ordered after each data refresh
$dataRange = $ws.Range("A4:H205")
$sortKey = $ws.Range("B4:B205")
$dataRange.Sort($sortKey, 2)
